$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet
$ws.Name = "Nodes"

# Insert a new header row at the top, shifting existing data down
$ws.Rows.Item(1).Insert()

# Set header values
$ws.Range("A1").Value = "Id"
$ws.Range("B1").Value = "SupplierName"
